# Refresh the crypto price/volume table to the latest scrape.
# Price (D) / Volume(1h) (E) columns hold text-formatted numbers
# (e.g. "1.00", "  -1.72%  "); Excel's COM Value setter would silently
# coerce plain-looking numerics to real numbers, so those cells are
# force-formatted as Text ("@") first to keep the literal string.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "61.823.97"
$ws.Range("E2").Value = "  -1.72%  "
$ws.Range("D3").Value = "3.409.36"
$ws.Range("E3").Value = "  -1.58%  "
$ws.Range("E4").Value = "  -0.03%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "404.61"
$ws.Range("E5").Value = "  -1.05%  "
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "131.72"
$ws.Range("E6").Value = "  +0.10%  "
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "0.590"
$ws.Range("E7").Value = "  -2.13%  "
$ws.Range("E8").Value = "  +0.04%  "
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "0.683"
$ws.Range("E9").Value = "  -1.36%  "
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "0.125"
$ws.Range("E10").Value = "  -3.40%  "
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "41.73"
$ws.Range("E11").Value = "  -2.91%  "
$ws.Range("E12").Value = "  -1.16%  "
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "8.40"
$ws.Range("E13").Value = "  -4.30%  "
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "19.74"
$ws.Range("E14").Value = "  -1.94%  "
$ws.Range("D15").Value = "3.434.84"
$ws.Range("E15").Value = "  -0.82%  "
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "11.58"
$ws.Range("E16").Value = "  +6.85%  "
$ws.Range("D17").Value = "61.839.72"
$ws.Range("E17").Value = "  -1.79%  "
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "1.01"
$ws.Range("E18").Value = "  -3.09%  "
$ws.Range("E19").Value = "  +2.29%  "
$ws.Range("E20").Value = "  -5.30%  "
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "83.33"
$ws.Range("E21").Value = "  +0.91%  "
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "311.04"
$ws.Range("E22").Value = "  -1.25%  "
$ws.Range("E23").Value = "  -2.44%  "
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "3.13"
$ws.Range("E24").Value = "  -0.79%  "
$ws.Range("E25").Value = "  +9.48%  "
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "29.51"
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "8.08"
$ws.Range("E27").Value = "  -1.55%  "
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "7.65"
$ws.Range("E28").Value = "  +0.38%  "
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "2.74"
$ws.Range("E29").Value = "  +4.93%  "
$ws.Range("E30").Value = "  -2.28%  "
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "0.114"
$ws.Range("E31").Value = "  -2.93%  "
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "42.33"
$ws.Range("E32").Value = "  -4.36%  "
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "1.00"
$ws.Range("E33").Value = "  +0.06%  "
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "11.31"
$ws.Range("E34").Value = "  -4.14%  "
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "0.0482"
$ws.Range("E35").Value = "  -2.09%  "
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "51.22"
$ws.Range("E36").Value = "  -2.52%  "
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "0.998"
$ws.Range("E37").Value = "  -0.13%  "
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "0.326"
$ws.Range("E38").Value = "  +13.70%  "
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "3.37"
$ws.Range("E39").Value = "  -5.64%  "
$ws.Range("E40").Value = "  -3.33%  "
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "138.97"
$ws.Range("E41").Value = "  +2.53%  "
$ws.Range("B42").Value = "Stellar"
$ws.Range("C42").Value = "https://coinranking.com/coin/f3iaFeCKEmkaZ+stellar-xlm"
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "0.124"
$ws.Range("E42").Value = "  -1.07%  "
$ws.Range("B43").Value = "ARBITRUM"
$ws.Range("C43").Value = "https://coinranking.com/coin/1Uo6s62Oc+arbitrum-arb"
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "1.97"
$ws.Range("E43").Value = "  -0.93%  "
$ws.Range("E44").Value = "  -0.24%  "
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "16.64"
$ws.Range("E45").Value = "  -4.18%  "
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "2.23"
$ws.Range("E46").Value = "  -1.09%  "
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "21.25"
$ws.Range("E47").Value = "  -3.72%  "
$ws.Range("D48").Value = "2.106.90"
$ws.Range("E48").Value = "  -3.49%  "
$ws.Range("E49").Value = "  -3.29%  "
$ws.Range("B50").Value = "ThetaToken"
$ws.Range("C50").Value = "https://coinranking.com/coin/B42IRxNtoYmwK+thetatoken-theta"
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "1.95"
$ws.Range("E50").Value = "  +4.19%  "
$ws.Range("B51").Value = "Fetch.AI"
$ws.Range("C51").Value = "https://coinranking.com/coin/AWma-WzFHmKVQ+fetchai-fet"
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "1.76"
$ws.Range("E51").Value = "  +20.96%  "
